# edit.ps1 — applies the LOB1205.docx restructuring described by the diff.
# Paragraph count/order/styles are unchanged by this edit; only the text content
# of specific paragraphs (and a run-level reshuffle inside the "Avaliação" bullet
# paragraph) moves around. We therefore edit by paragraph index, and for the one
# paragraph whose internal runs get reordered we use scoped Find/Replace so the
# existing bold "label" runs (Método:/Critério:/Norma de recuperação:) keep their
# formatting untouched.

$d = $word.ActiveDocument
$vtab = [char]11   # w:br <-> vertical-tab char in Range.Text / Find

# --- Paragraph 6: "Objetivos" PT paragraph becomes the PT numbered program list ---
$d.Paragraphs.Item(6).Range.Text = '1. Problemas ambientais, causas e soluções' + $vtab + '2. Ecossistemas: o que são e como funcionam' + $vtab + '3. Biodiversidade e evolução' + $vtab + '4. Biodiversidade, interações de espécies e controle da população' + $vtab + '5. A população humana e seu impacto'

# --- Paragraph 7: "Objetivos" EN (italic) paragraph becomes the EN numbered list ---
$d.Paragraphs.Item(7).Range.Text = '1. Environmental problems, causes and solutions' + $vtab + '2. Ecosystems: what they are and how they work' + $vtab + '3. Biodiversity and evolution' + $vtab + '4. Biodiversity, species interactions and population control' + $vtab + '5. The human population and its impact'

# --- Paragraph 9: docente bullet becomes the PT objectives paragraph ---
$d.Paragraphs.Item(9).Range.Text = 'Abordar os princípios e conceitos da evolução biológica e da ecologia em suas diferentes escalas: populações, comunidades e ecossistemas, dentro do enfoque da sustentabilidade. Apresentar aspectos econômicos, sociais e culturais da sociedade envolvidos na preservação da biodiversidade e dos ecossistemas.'

# --- Paragraph 12: EN numbered list (italic) becomes the EN objectives paragraph ---
$d.Paragraphs.Item(12).Range.Text = 'To approach the principles and concepts of biological evolution and ecology in its different scales: populations, communities and ecosystems, within the focus of sustainability. To present economic, social and cultural aspects of society involved in the preservation of biodiversity and ecosystems.'

# --- Paragraph 14: PT numbered list becomes the "Método" evaluation text ---
$d.Paragraphs.Item(14).Range.Text = 'A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e seminários.'

# --- Paragraph 19: bibliography paragraph becomes the docente bullet text ---
$d.Paragraphs.Item(19).Range.Text = '5840942 - Marco Aurélio Kondracki de Alcântara'

# --- Paragraph 17 ("Avaliação" bullet list): the text that follows each bold
#     label ("Método: "/"Critério: "/"Norma de recuperação: ") shifts down one
#     slot, and the final slot becomes the bibliography text that used to live in
#     its own "Bibliografia" paragraph. The bold label runs themselves are left
#     untouched; we locate-and-replace the plain runs after them in document order
#     so each scoped Find starts after the previous replacement.
$p17 = $d.Paragraphs.Item(17)
$pStart = $p17.Range.Start
$pEnd = $p17.Range.End

# Step 1: text after "Método: " (was the P1/P2 method blurb) -> becomes the
#         "O aluno poderá optar..." criteria text (with its internal breaks).
$r1 = $d.Range($pStart, $pEnd)
$target1 = 'A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e seminários.'
$replacement1 = 'O aluno poderá optar por dois critérios de avaliação:' + $vtab + 'Critério 1: NF = (P1+P2)/2; ou' + $vtab + 'Critério 2: NF = (NOTA 1 + NOTA 2)/2' + $vtab + 'Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas.'
$null = $r1.Find.Execute($target1, $false, $false, $false, $false, $false, $true, 1, $false, $replacement1, 2)

# Step 2: text after "Critério: " (was the criteria text) -> becomes the
#         "Exame Final..." recovery-norm text. Search starts after step 1's match
#         so the (identical-looking) freshly-inserted copy is not matched again.
$pEnd2 = $d.Paragraphs.Item(17).Range.End
$r2 = $d.Range($r1.End, $pEnd2)
$target2 = 'O aluno poderá optar por dois critérios de avaliação:' + $vtab + 'Critério 1: NF = (P1+P2)/2; ou' + $vtab + 'Critério 2: NF = (NOTA 1 + NOTA 2)/2' + $vtab + 'Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas.'
$replacement2 = 'Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.'
$null = $r2.Find.Execute($target2, $false, $false, $false, $false, $false, $true, 1, $false, $replacement2, 2)

# Step 3: text after "Norma de recuperação: " (was the recovery-norm text) ->
#         becomes the full bibliography body (previously its own "Bibliografia" paragraph).
$pEnd3 = $d.Paragraphs.Item(17).Range.End
$r3 = $d.Range($r2.End, $pEnd3)
$target3 = 'Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.'
$replacement3 = 'Básica:' + $vtab + 'MILLER, G.T.; SPOOLMAN, S.E. 2012. Ecologia e sustentabilidade. Cengage Learning. 412p.' + $vtab + '' + $vtab + 'Complementar:' + $vtab + 'BEGON, M., J.L. HARPER & C.R. TOWNSEND. 2005. Ecology. From Individuals to Communities. Blackwell Science.' + $vtab + 'RICKLEFS, R.E. 2003. A economia da natureza. Guanabara Koogan.' + $vtab + 'RICKLEFS, R.E. & G.L. MILLER. 2000. Ecology. W.H. Freeman and Co.' + $vtab + 'TOWNSEND, C.R., M. BEGON. & J.L. HARPER 2006. Fundamentos em ecologia. Artmed.'
$null = $r3.Find.Execute($target3, $false, $false, $false, $false, $false, $true, 1, $false, $replacement3, 2)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
